# "Add per user signaling and linking"
#
# - lifts: rows 3-5 move from STATE=NONE to STATE=READY; row 4's NOTE becomes
#   "🅱atteri"; row 5's NOTE becomes the text "420"; a new USERS column (F) is
#   populated for rows 2-5; the old rows 6-7 (OPENING/READY placeholder rows)
#   are removed.
# - follows: the site codes that used to live across B1:C1 are turned into
#   their own rows in column A (A2, A3), and B1:C1 are cleared out.

$wb = $excel.ActiveWorkbook

# ---- Sheet "follows" (do this first so the final active tab ends up being "lifts") ----
$ws2 = $wb.Worksheets.Item("follows")
$ws2.Activate()

$site1 = $ws2.Range("B1").Value2
$site2 = $ws2.Range("C1").Value2

$ws2.Columns.Item(1).Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

$ws2.Range("A2").Value = $site1
$ws2.Range("A3").Value = $site2

$ws2.Range("B1:C1").Clear()

$ws2.Range("B1").Select()

# ---- Sheet "lifts" ----
$ws = $wb.Worksheets.Item("lifts")
$ws.Activate()

# Update STATE column (NONE -> READY) for rows 3-5
$ws.Range("B3").Value = "READY"
$ws.Range("B4").Value = "READY"
$ws.Range("B5").Value = "READY"

# Update NOTE column for rows 4-5
$ws.Range("E4").Value = "🅱atteri"

# E5 becomes the text "420" (not the number 420) - round-trip through TEXT()
# + paste-values so it lands as a plain string without picking up a new
# number-format style.
$ws.Range("E5").Formula = '=TEXT(420,"0")'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# New USERS column (F) values
$ws.Range("F2").Value = "Eemeli"
$ws.Range("F3").Value = "Eemelioma"
$ws.Range("F4").Value = "Eemeli:Eemelioma:Akseli"
$ws.Range("F5").Value = "Eemelioma:Akseli"

# F3:F5 keep the default (unstyled) look; only F2 is centered like the rest of the table
$ws.Range("F3:F5").Style = "Normal"

# Remove the old rows 6 and 7 (OPENING / READY placeholder rows)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

$ws.Range("F10").Select()
